$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from BA1 into the new header cells
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BA1").Copy($ws.Range("BC1"))
$ws.Range("BA1").Copy($ws.Range("BD1"))

$ws.Range("BB1").Value = "31/12/2023"
$ws.Range("BC1").Value = "31/03/2024"
$ws.Range("BD1").Value = "30/06/2024"

$ws.Range("BB2").Value = 665614.976 ; $ws.Range("BC2").Value = 639088 ; $ws.Range("BD2").Value = 636880
$ws.Range("BB3").Value = 437689.984 ; $ws.Range("BC3").Value = 412012.992 ; $ws.Range("BD3").Value = 407015.008
$ws.Range("BB4").Value = 125152 ; $ws.Range("BC4").Value = 112328 ; $ws.Range("BD4").Value = 93939
$ws.Range("BB5").Value = 119 ; $ws.Range("BC5").Value = 97 ; $ws.Range("BD5").Value = 322
$ws.Range("BB6").Value = 152111.008 ; $ws.Range("BC6").Value = 136135.008 ; $ws.Range("BD6").Value = 154034
$ws.Range("BB7").Value = 117525 ; $ws.Range("BC7").Value = 121306 ; $ws.Range("BD7").Value = 112982
$ws.Range("BB8").Value = 0 ; $ws.Range("BC8").Value = 0 ; $ws.Range("BD8").Value = 0
$ws.Range("BB9").Value = 26271 ; $ws.Range("BC9").Value = 23708 ; $ws.Range("BD9").Value = 24950
$ws.Range("BB10").Value = 0 ; $ws.Range("BC10").Value = 0 ; $ws.Range("BD10").Value = 0
$ws.Range("BB11").Value = 16512 ; $ws.Range("BC11").Value = 18439 ; $ws.Range("BD11").Value = 20788
$ws.Range("BB12").Value = 9328 ; $ws.Range("BC12").Value = 9062 ; $ws.Range("BD12").Value = 10763
$ws.Range("BB13").Value = 0 ; $ws.Range("BC13").Value = 0 ; $ws.Range("BD13").Value = 0
$ws.Range("BB14").Value = 0 ; $ws.Range("BC14").Value = 0 ; $ws.Range("BD14").Value = 0
$ws.Range("BB15").Value = 0 ; $ws.Range("BC15").Value = 0 ; $ws.Range("BD15").Value = 0
$ws.Range("BB16").Value = 6076 ; $ws.Range("BC16").Value = 5810 ; $ws.Range("BD16").Value = 7511
$ws.Range("BB17").Value = 0 ; $ws.Range("BC17").Value = 0 ; $ws.Range("BD17").Value = 0
$ws.Range("BB18").Value = 0 ; $ws.Range("BC18").Value = 0 ; $ws.Range("BD18").Value = 0
$ws.Range("BB19").Value = 0 ; $ws.Range("BC19").Value = 0 ; $ws.Range("BD19").Value = 0
$ws.Range("BB20").Value = 0 ; $ws.Range("BC20").Value = 0 ; $ws.Range("BD20").Value = 0
$ws.Range("BB21").Value = 0 ; $ws.Range("BC21").Value = 0 ; $ws.Range("BD21").Value = 0
$ws.Range("BB22").Value = 0 ; $ws.Range("BC22").Value = 0 ; $ws.Range("BD22").Value = 0
$ws.Range("BB23").Value = 27556 ; $ws.Range("BC23").Value = 26963 ; $ws.Range("BD23").Value = 27479
$ws.Range("BB24").Value = 191040.992 ; $ws.Range("BC24").Value = 191050 ; $ws.Range("BD24").Value = 191623.008
$ws.Range("BB25").Value = 0 ; $ws.Range("BC25").Value = 0 ; $ws.Range("BD25").Value = 0
$ws.Range("BB26").Value = 665614.976 ; $ws.Range("BC26").Value = 639088 ; $ws.Range("BD26").Value = 636880
$ws.Range("BB27").Value = 95547 ; $ws.Range("BC27").Value = 83353 ; $ws.Range("BD27").Value = 79794
$ws.Range("BB28").Value = 10780 ; $ws.Range("BC28").Value = 8041 ; $ws.Range("BD28").Value = 10505
$ws.Range("BB29").Value = 44931 ; $ws.Range("BC29").Value = 41044 ; $ws.Range("BD29").Value = 36724
$ws.Range("BB30").Value = 8457 ; $ws.Range("BC30").Value = 6121 ; $ws.Range("BD30").Value = 6061
$ws.Range("BB31").Value = 16060 ; $ws.Range("BC31").Value = 17536 ; $ws.Range("BD31").Value = 16980
$ws.Range("BB32").Value = 0 ; $ws.Range("BC32").Value = 0 ; $ws.Range("BD32").Value = 0
$ws.Range("BB33").Value = 0 ; $ws.Range("BC33").Value = 0 ; $ws.Range("BD33").Value = 0
$ws.Range("BB34").Value = 15319 ; $ws.Range("BC34").Value = 10362 ; $ws.Range("BD34").Value = 9524
$ws.Range("BB35").Value = 0 ; $ws.Range("BC35").Value = 249 ; $ws.Range("BD35").Value = 0
$ws.Range("BB36").Value = 0 ; $ws.Range("BC36").Value = 0 ; $ws.Range("BD36").Value = 0
$ws.Range("BB37").Value = 167388.992 ; $ws.Range("BC37").Value = 159103.008 ; $ws.Range("BD37").Value = 155212.992
$ws.Range("BB38").Value = 77759 ; $ws.Range("BC38").Value = 72598 ; $ws.Range("BD38").Value = 69021
$ws.Range("BB39").Value = 0 ; $ws.Range("BC39").Value = 0 ; $ws.Range("BD39").Value = 0
$ws.Range("BB40").Value = 1408 ; $ws.Range("BC40").Value = 1408 ; $ws.Range("BD40").Value = 1311
$ws.Range("BB41").Value = 26813 ; $ws.Range("BC41").Value = 29618 ; $ws.Range("BD41").Value = 33286
$ws.Range("BB42").Value = 0 ; $ws.Range("BC42").Value = 0 ; $ws.Range("BD42").Value = 0
$ws.Range("BB43").Value = 61409 ; $ws.Range("BC43").Value = 55479 ; $ws.Range("BD43").Value = 51595
$ws.Range("BB44").Value = 0 ; $ws.Range("BC44").Value = 0 ; $ws.Range("BD44").Value = 0
$ws.Range("BB45").Value = 0 ; $ws.Range("BC45").Value = 0 ; $ws.Range("BD45").Value = 0
$ws.Range("BB46").Value = 0 ; $ws.Range("BC46").Value = 0 ; $ws.Range("BD46").Value = 0
$ws.Range("BB47").Value = 402679.008 ; $ws.Range("BC47").Value = 396632 ; $ws.Range("BD47").Value = 401872.992
$ws.Range("BB48").Value = 130583 ; $ws.Range("BC48").Value = 130583 ; $ws.Range("BD48").Value = 130583
$ws.Range("BB49").Value = 172784.992 ; $ws.Range("BC49").Value = 170712.992 ; $ws.Range("BD49").Value = 160202
$ws.Range("BB50").Value = 0 ; $ws.Range("BC50").Value = 0 ; $ws.Range("BD50").Value = 0
$ws.Range("BB51").Value = 124039 ; $ws.Range("BC51").Value = 114300 ; $ws.Range("BD51").Value = 110352
$ws.Range("BB52").Value = 0 ; $ws.Range("BC52").Value = 5871 ; $ws.Range("BD52").Value = 25941
$ws.Range("BB53").Value = -13858 ; $ws.Range("BC53").Value = -13965 ; $ws.Range("BD53").Value = -14335
$ws.Range("BB54").Value = 0 ; $ws.Range("BC54").Value = 0 ; $ws.Range("BD54").Value = 0
$ws.Range("BB55").Value = -10870 ; $ws.Range("BC55").Value = -10870 ; $ws.Range("BD55").Value = -10870
$ws.Range("BB56").Value = 0 ; $ws.Range("BC56").Value = 0 ; $ws.Range("BD56").Value = 0
# Row 57 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA57").Copy($ws.Range("BB57"))
$ws.Range("BA57").Copy($ws.Range("BC57"))
$ws.Range("BA57").Copy($ws.Range("BD57"))

# Row 58 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA58").Copy($ws.Range("BB58"))
$ws.Range("BA58").Copy($ws.Range("BC58"))
$ws.Range("BA58").Copy($ws.Range("BD58"))

$ws.Range("BB59").Value = 111067.008 ; $ws.Range("BC59").Value = 68766 ; $ws.Range("BD59").Value = 99556
$ws.Range("BB60").Value = -48579 ; $ws.Range("BC60").Value = -32256 ; $ws.Range("BD60").Value = -44816
$ws.Range("BB61").Value = 62487.992 ; $ws.Range("BC61").Value = 36510 ; $ws.Range("BD61").Value = 54740
$ws.Range("BB62").Value = -23583 ; $ws.Range("BC62").Value = -21942 ; $ws.Range("BD62").Value = -25880
$ws.Range("BB63").Value = -10640 ; $ws.Range("BC63").Value = -9989 ; $ws.Range("BD63").Value = -10547
$ws.Range("BB64").Value = -1630 ; $ws.Range("BC64").Value = 0 ; $ws.Range("BD64").Value = 0
$ws.Range("BB65").Value = 0 ; $ws.Range("BC65").Value = 0 ; $ws.Range("BD65").Value = 0
$ws.Range("BB66").Value = -574 ; $ws.Range("BC66").Value = -235 ; $ws.Range("BD66").Value = 823
$ws.Range("BB67").Value = 0 ; $ws.Range("BC67").Value = 0 ; $ws.Range("BD67").Value = 0
$ws.Range("BB68").Value = 3300 ; $ws.Range("BC68").Value = 5056 ; $ws.Range("BD68").Value = 6329
$ws.Range("BB69").Value = 7301 ; $ws.Range("BC69").Value = 6423 ; $ws.Range("BD69").Value = 9780
$ws.Range("BB70").Value = -4001 ; $ws.Range("BC70").Value = -1367 ; $ws.Range("BD70").Value = -3451
# Row 71 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA71").Copy($ws.Range("BB71"))
$ws.Range("BA71").Copy($ws.Range("BC71"))
$ws.Range("BA71").Copy($ws.Range("BD71"))

# Row 72 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA72").Copy($ws.Range("BB72"))
$ws.Range("BA72").Copy($ws.Range("BC72"))
$ws.Range("BA72").Copy($ws.Range("BD72"))

# Row 73 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA73").Copy($ws.Range("BB73"))
$ws.Range("BA73").Copy($ws.Range("BC73"))
$ws.Range("BA73").Copy($ws.Range("BD73"))

$ws.Range("BB74").Value = 29361 ; $ws.Range("BC74").Value = 9400 ; $ws.Range("BD74").Value = 25465
$ws.Range("BB75").Value = 103 ; $ws.Range("BC75").Value = -724 ; $ws.Range("BD75").Value = -1727
$ws.Range("BB76").Value = -4340 ; $ws.Range("BC76").Value = -2805 ; $ws.Range("BD76").Value = -3668
# Row 77 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA77").Copy($ws.Range("BB77"))
$ws.Range("BA77").Copy($ws.Range("BC77"))
$ws.Range("BA77").Copy($ws.Range("BD77"))

# Row 78 is blank; replicate the blank placeholder cells from column BA
$ws.Range("BA78").Copy($ws.Range("BB78"))
$ws.Range("BA78").Copy($ws.Range("BC78"))
$ws.Range("BA78").Copy($ws.Range("BD78"))

$ws.Range("BB79").Value = 0 ; $ws.Range("BC79").Value = 0 ; $ws.Range("BD79").Value = 0
$ws.Range("BB80").Value = 25124 ; $ws.Range("BC80").Value = 5871 ; $ws.Range("BD80").Value = 20070
